# Auto-generated Excel COM-interop script
# Updates cached market-price / leve-profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4513
$ws.Range("I62").Value = 4462.25
$ws.Range("J62").Value = 4571
$ws.Range("K62").Value = 4462.25
$ws.Range("L62").Value = 4571
$ws.Range("M62").Value = -3838.25
$ws.Range("N62").Value = -5819

$ws.Range("H65").Value = 4513
$ws.Range("I65").Value = 4462.25
$ws.Range("J65").Value = 4571
$ws.Range("K65").Value = 22311.25
$ws.Range("L65").Value = 22855
$ws.Range("M65").Value = -19191.25
$ws.Range("N65").Value = -29095

$ws.Range("H76").Value = 3999
$ws.Range("I76").Value = 1999
$ws.Range("K76").Value = 1999
$ws.Range("M76").Value = -1684

$ws.Range("H79").Value = 3999
$ws.Range("I79").Value = 1999
$ws.Range("K79").Value = 1999
$ws.Range("M79").Value = -907

$ws.Range("H140").Value = 66136.586
$ws.Range("J140").Value = 65720.91
$ws.Range("L140").Value = 65720.91
$ws.Range("N140").Value = -76080.91

$ws.Range("H141").Value = 9372.625
$ws.Range("I141").Value = 8022.913
$ws.Range("J141").Value = 12821.889
$ws.Range("K141").Value = 24068.739
$ws.Range("L141").Value = 38465.667
$ws.Range("M141").Value = -18888.739
$ws.Range("N141").Value = -48825.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4526.976
$ws.Range("I61").Value = 3680.7932
$ws.Range("J61").Value = 6414.615
$ws.Range("K61").Value = 3680.7932
$ws.Range("L61").Value = 6414.615
$ws.Range("M61").Value = -3468.7932
$ws.Range("N61").Value = -6838.615

$ws.Range("H92").Value = 30000000
$ws.Range("J92").Value = 30000000
$ws.Range("L92").Value = 30000000
$ws.Range("N92").Value = -30004992

$ws.Range("H132").Value = 21900.941
$ws.Range("I132").Value = 29700.9
$ws.Range("J132").Value = 10758.143
$ws.Range("K132").Value = 89102.70000000001
$ws.Range("L132").Value = 32274.429
$ws.Range("M132").Value = -86572.70000000001
$ws.Range("N132").Value = -37334.429

$ws.Range("H136").Value = 4526.976
$ws.Range("I136").Value = 3680.7932
$ws.Range("J136").Value = 6414.615
$ws.Range("K136").Value = 11042.3796
$ws.Range("L136").Value = 19243.845
$ws.Range("M136").Value = -8492.3796
$ws.Range("N136").Value = -24343.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 46.363636
$ws.Range("I11").Value = 50.9
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 50.9
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 89.09999999999999
$ws.Range("N11").Value = -281

$ws.Range("H94").Value = 623722.25
$ws.Range("I94").Value = 857069.1
$ws.Range("J94").Value = 1464
$ws.Range("K94").Value = 857069.1
$ws.Range("L94").Value = 1464
$ws.Range("M94").Value = -856618.1
$ws.Range("N94").Value = -2366

$ws.Range("H99").Value = 1043213.2
$ws.Range("I99").Value = 1226903.1
$ws.Range("K99").Value = 1226903.1
$ws.Range("M99").Value = -1225405.1

$ws.Range("H134").Value = 2234.55
$ws.Range("I134").Value = 1648.7222
$ws.Range("J134").Value = 7507
$ws.Range("K134").Value = 4946.1666
$ws.Range("L134").Value = 22521
$ws.Range("M134").Value = -2411.1666
$ws.Range("N134").Value = -27591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12822583
$ws.Range("I31").Value = 14707492
$ws.Range("J31").Value = 5195.2
$ws.Range("K31").Value = 14707492
$ws.Range("L31").Value = 5195.2
$ws.Range("M31").Value = -14707197
$ws.Range("N31").Value = -5785.2

$ws.Range("H33").Value = 3004
$ws.Range("I33").Value = 3004
$ws.Range("K33").Value = 3004
$ws.Range("M33").Value = -2625

$ws.Range("H34").Value = 12822583
$ws.Range("I34").Value = 14707492
$ws.Range("J34").Value = 5195.2
$ws.Range("K34").Value = 14707492
$ws.Range("L34").Value = 5195.2
$ws.Range("M34").Value = -14707290
$ws.Range("N34").Value = -5599.2

$ws.Range("H132").Value = 95247830
$ws.Range("I132").Value = 111113304
$ws.Range("J132").Value = 55000
$ws.Range("K132").Value = 333339912
$ws.Range("L132").Value = 165000
$ws.Range("M132").Value = -333337382
$ws.Range("N132").Value = -170060

$ws.Range("H134").Value = 2808.0417
$ws.Range("I134").Value = 2528.2354
$ws.Range("K134").Value = 7584.706200000001
$ws.Range("M134").Value = -5049.706200000001

$ws.Range("H141").Value = 108696.68
$ws.Range("J141").Value = 109193.98
$ws.Range("L141").Value = 109193.98
$ws.Range("N141").Value = -119553.98

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 6221.875
$ws.Range("J52").Value = 6221.875
$ws.Range("L52").Value = 18665.625
$ws.Range("N52").Value = -19197.625

$ws.Range("H109").Value = 11000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 11000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 33000
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -35080

$ws.Range("H114").Value = 2038.375
$ws.Range("I114").Value = 839.6667
$ws.Range("J114").Value = 2315
$ws.Range("K114").Value = 2519.0001
$ws.Range("L114").Value = 6945
$ws.Range("M114").Value = 734.9998999999998
$ws.Range("N114").Value = -13453

$ws.Range("H117").Value = 3036.25
$ws.Range("J117").Value = 2902.4546
$ws.Range("L117").Value = 8707.363799999999
$ws.Range("N117").Value = -15591.3638

$ws.Range("H134").Value = 4287
$ws.Range("I134").Value = 1491.9
$ws.Range("J134").Value = 11274.75
$ws.Range("K134").Value = 4475.700000000001
$ws.Range("L134").Value = 33824.25
$ws.Range("M134").Value = 594.2999999999993
$ws.Range("N134").Value = -43964.25

$ws.Range("H139").Value = 2904.9167
$ws.Range("I139").Value = 809.8333
$ws.Range("K139").Value = 2429.4999
$ws.Range("M139").Value = 2710.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 242.6
$ws.Range("I31").Value = 295.75
$ws.Range("J31").Value = 30
$ws.Range("K31").Value = 295.75
$ws.Range("L31").Value = 30
$ws.Range("M31").Value = -3.75
$ws.Range("N31").Value = -614

$ws.Range("H37").Value = 242.6
$ws.Range("I37").Value = 295.75
$ws.Range("J37").Value = 30
$ws.Range("K37").Value = 295.75
$ws.Range("L37").Value = 30
$ws.Range("M37").Value = -18.75
$ws.Range("N37").Value = -584

$ws.Range("H97").Value = 743.25
$ws.Range("I97").Value = 655.8333
$ws.Range("J97").Value = 1005.5
$ws.Range("K97").Value = 655.8333
$ws.Range("L97").Value = 1005.5
$ws.Range("M97").Value = -159.8333
$ws.Range("N97").Value = -1997.5

$ws.Range("H107").Value = 6803590.5
$ws.Range("I107").Value = 11905213
$ws.Range("K107").Value = 11905213
$ws.Range("M107").Value = -11903293

$ws.Range("H113").Value = 997
$ws.Range("I113").Value = 997
$ws.Range("K113").Value = 997
$ws.Range("M113").Value = 1173

$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884

$ws.Range("H132").Value = 3275.4
$ws.Range("I132").Value = 2795.0715
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 8385.2145
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -5855.2145
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4035.275
$ws.Range("J7").Value = 5402.32
$ws.Range("L7").Value = 5402.32
$ws.Range("N7").Value = -5626.32

$ws.Range("H93").Value = 1452.0416
$ws.Range("I93").Value = 1332.0667
$ws.Range("K93").Value = 1332.0667
$ws.Range("M93").Value = -84.06670000000008

$ws.Range("H104").Value = 43706.11
$ws.Range("J104").Value = 43706.11
$ws.Range("L104").Value = 43706.11
$ws.Range("N104").Value = -50694.11

$ws.Range("H118").Value = 120000
$ws.Range("J118").Value = 120000
$ws.Range("L118").Value = 120000
$ws.Range("N118").Value = -123314

$ws.Range("H122").Value = 8123.276
$ws.Range("I122").Value = 3908.3333
$ws.Range("J122").Value = 12639.286
$ws.Range("K122").Value = 11724.9999
$ws.Range("L122").Value = 37917.858
$ws.Range("M122").Value = -9274.999899999999
$ws.Range("N122").Value = -42817.858

$ws.Range("H126").Value = 4035.275
$ws.Range("J126").Value = 5402.32
$ws.Range("L126").Value = 16206.96
$ws.Range("N126").Value = -21146.96

$ws.Range("H132").Value = 2661.12
$ws.Range("I132").Value = 2614.6494
$ws.Range("J132").Value = 2816.6956
$ws.Range("K132").Value = 7843.948199999999
$ws.Range("L132").Value = 8450.086800000001
$ws.Range("M132").Value = -5313.948199999999
$ws.Range("N132").Value = -13510.0868

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 49974
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 49974
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 49974
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -50434

$ws.Range("H62").Value = 7666.1665
$ws.Range("J62").Value = 9499.375
$ws.Range("L62").Value = 9499.375
$ws.Range("N62").Value = -10747.375

$ws.Range("H65").Value = 7666.1665
$ws.Range("J65").Value = 9499.375
$ws.Range("L65").Value = 47496.875
$ws.Range("N65").Value = -53736.875

$ws.Range("H82").Value = 10000
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 10000
$ws.Range("N82").Value = -10766

$ws.Range("H85").Value = 10000
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 10000
$ws.Range("N85").Value = -12652

$ws.Range("H96").Value = 5000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -12746

